$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F "想去人数" (people who want to go) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1329
$ws1.Range("F4").Value = 1133
$ws1.Range("F5").Value = 1023
$ws1.Range("F6").Value = 1797
$ws1.Range("F7").Value = 565
$ws1.Range("F8").Value = 1205
$ws1.Range("F9").Value = 57
$ws1.Range("F12").Value = 298
$ws1.Range("F13").Value = 69
$ws1.Range("F15").Value = 691
$ws1.Range("F16").Value = 171
$ws1.Range("F17").Value = 104
$ws1.Range("F20").Value = 329
$ws1.Range("F21").Value = 151
$ws1.Range("F22").Value = 670
$ws1.Range("F23").Value = 35
$ws1.Range("F24").Value = 647
$ws1.Range("F25").Value = 150
$ws1.Range("F27").Value = 876
$ws1.Range("F28").Value = 316
$ws1.Range("F31").Value = 272

# Sheet "演出" (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 320
$ws2.Range("F7").Value = 255

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1329
$ws4.Range("F5").Value = 1133
$ws4.Range("F6").Value = 1023
$ws4.Range("F7").Value = 1797
$ws4.Range("F8").Value = 565
$ws4.Range("F9").Value = 1205
$ws4.Range("F10").Value = 57
$ws4.Range("F14").Value = 298
$ws4.Range("F15").Value = 69
$ws4.Range("F17").Value = 691
$ws4.Range("F18").Value = 171
$ws4.Range("F19").Value = 104
$ws4.Range("F22").Value = 320
$ws4.Range("F25").Value = 329
$ws4.Range("F27").Value = 255
$ws4.Range("F28").Value = 255
$ws4.Range("F29").Value = 152
$ws4.Range("F30").Value = 670
$ws4.Range("F31").Value = 35
$ws4.Range("F32").Value = 647
$ws4.Range("F33").Value = 150
$ws4.Range("F35").Value = 876
$ws4.Range("F36").Value = 316
$ws4.Range("F41").Value = 272
